$wb = $excel.ActiveWorkbook

$namesSheet = $wb.Worksheets.Item("Sheet1")
$usedSheet = $wb.Worksheets.Item("used")

# The id currently sitting at the top of the pending "names" queue (row 2,
# since row 1 is the already-consumed/header anchor).
$id = $namesSheet.Range("A2").Value2

# Remove that row from the names sheet; everything below shifts up.
$namesSheet.Rows.Item(2).Delete()

# Find the first empty row in the "used" log and append the record there.
$lastRow = $usedSheet.Cells.Item($usedSheet.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$usedSheet.Cells.Item($newRow, 1).Value = $id
$usedSheet.Cells.Item($newRow, 2).Value = "ChatGPT Image 2026年1月21日 21_36_51.png"
$usedSheet.Cells.Item($newRow, 3).Value = "2026-01-21 21:37:45"
